$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers in M1:O1
$ws.Range("M1").Value = "renewd"
$ws.Range("N1").Value = "PlanID"
$ws.Range("O1").Value = "iteration"

# Copy format from L1 (existing header) to M1:O1
$ws.Range("L1").Copy()
$ws.Range("M1:O1").PasteSpecial(-4122)

# Fill data rows 2-31 for columns M, N, O
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 13).Value = "after"
    $ws.Cells.Item($r, 14).Value = 20110270
    $ws.Cells.Item($r, 15).Value = 12
}
